$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
Write-Host "ThemeColorScheme=" $tcs
Write-Host "Count=" $tcs.Count
$tcs | Get-Member | ForEach-Object { Write-Host $_.Name $_.MemberType $_.Definition }
